$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Visit" column (D) is stored as text in the workbook, so force a text
# number format on the rows we touch before writing the values, keeping them
# as strings (not converted to numbers).
$ws.Range("D2:D8").NumberFormat = "@"

# New values for rows 2,3,6,7,8 (a cyclic re-shuffle of the existing rows,
# likely caused by re-sorting in R before an incomplete save).
$ws.Range("A2").Value = "Participant-Visit Charlie-4 was not found in the Visit Enrollment sheet."
$ws.Range("B2").Value = "Visit Enrollment"
$ws.Range("C2").Value = "Charlie"
$ws.Range("D2").Value = "4"

$ws.Range("A3").Value = "Participant-Visit Charlie-4 was not found in LV."
$ws.Range("B3").Value = "LabVantage"
$ws.Range("C3").Value = "Charlie"
$ws.Range("D3").Value = "4"

$ws.Range("A6").Value = "Participant-Visit Alpha-2 was not found in LV."
$ws.Range("B6").Value = "LabVantage"
$ws.Range("C6").Value = "Alpha"
$ws.Range("D6").Value = "2"

$ws.Range("A7").Value = "Participant-Visit Bravo-2 was not found in the Visit Enrollment sheet."
$ws.Range("B7").Value = "Visit Enrollment"
$ws.Range("C7").Value = "Bravo"
$ws.Range("D7").Value = "2"

$ws.Range("A8").Value = "Participant-Visit Bravo-1 was not found in the BioBank sheet."
$ws.Range("B8").Value = "BioBank"
$ws.Range("C8").Value = "Bravo"
$ws.Range("D8").Value = "1"

$wb.Save()
